# Added prevalence over time into spreadsheet
$wb = $excel.ActiveWorkbook

# --- 1. Insert the new "prevelence_by_year" sheet, right before "fact" ---
$factSheet = $wb.Worksheets.Item("fact")
$newSheet = $wb.Worksheets.Add($factSheet)
$newSheet.Name = "prevelence_by_year"

# --- 2. Rename "summarymatrix" -> "summary_matrix" ---
$wb.Worksheets.Item("summarymatrix").Name = "summary_matrix"

# --- 3. Header row for the new sheet (bold, no border - matches style xf 101) ---
$newSheet.Cells.Item(1,1).Value = "code"
$newSheet.Cells.Item(1,2).Value = "year"
$newSheet.Cells.Item(1,3).Value = "t1dm"
$newSheet.Cells.Item(1,4).Value = "t2dm"
$newSheet.Cells.Item(1,5).Value = "other"
$newSheet.Cells.Item(1,6).Value = "midyrpe"
$newSheet.Range("A1:F1").Font.Bold = $true

# --- 4. Data rows (board code, year, t1dm, t2dm, other, midyrpe) ---
$rows = @(
@("S08000001",2011,2221,17919,35,366860),
@("S08000002",2011,614,4846,33,112870),
@("S08000003",2011,893,7236,39,148190),
@("S08000004",2011,1969,16164,59,364945),
@("S08000005",2011,1606,12528,67,293386),
@("S08000006",2011,3053,20902,85,550620),
@("S08000007",2011,6180,50005,527,1203870),
@("S08000008",2011,1758,12479,128,310830),
@("S08000009",2011,3513,24998,118,562477),
@("S08000010",2011,4175,29551,298,836711),
@("S08000011",2011,120,853,0,20110),
@("S08000012",2011,124,871,4,22400),
@("S08000013",2011,1864,18104,98,402641),
@("S08000014",2011,182,1058,1,26190),
@("S08000001",2010,2238,16775,62,367160),
@("S08000002",2010,601,4728,26,112680),
@("S08000003",2010,888,6836,47,148510),
@("S08000004",2010,1911,15480,76,363385),
@("S08000005",2010,1568,12007,43,291383),
@("S08000006",2010,3045,20227,85,544980),
@("S08000007",2010,6115,48090,265,1199026),
@("S08000008",2010,1706,12100,108,310530),
@("S08000009",2010,3480,23840,130,562215),
@("S08000010",2010,4109,28279,329,826231),
@("S08000011",2010,116,807,0,19960),
@("S08000012",2010,119,834,5,22210),
@("S08000013",2010,1837,17283,103,399550),
@("S08000014",2010,177,993,0,26180),
@("S08000001",2009,2234,15754,21,367510),
@("S08000002",2009,596,4530,5,112430),
@("S08000003",2009,871,6453,16,148580),
@("S08000004",2009,1896,14718,52,361815),
@("S08000005",2009,1526,11543,11,290047),
@("S08000006",2009,2976,19361,53,539630),
@("S08000007",2009,5923,46345,203,1194675),
@("S08000008",2009,1688,11470,28,309900),
@("S08000009",2009,3454,22794,54,561174),
@("S08000010",2009,4019,27506,215,817727),
@("S08000011",2009,118,776,0,19890),
@("S08000012",2009,114,792,1,21980),
@("S08000013",2009,1771,16283,42,396942),
@("S08000014",2009,181,939,0,26200),
@("S08000001",2008,2209,14764,19,367020),
@("S08000002",2008,596,4295,11,111430),
@("S08000003",2008,884,6098,21,148300),
@("S08000004",2008,1826,14041,70,360428),
@("S08000005",2008,1501,11153,7,288473),
@("S08000006",2008,2971,18048,126,535290),
@("S08000007",2008,6348,45639,346,1192419),
@("S08000008",2008,1673,11131,38,308790),
@("S08000009",2008,3415,21452,48,560042),
@("S08000010",2008,3933,26240,221,809764),
@("S08000011",2008,119,746,0,19860),
@("S08000012",2008,115,752,1,21950),
@("S08000013",2008,1692,15530,47,394134),
@("S08000014",2008,182,883,1,26300),
@("S08000001",2007,2235,13590,87,366450),
@("S08000002",2007,578,3972,46,110247),
@("S08000003",2007,913,5745,78,148030),
@("S08000004",2007,1860,13291,220,358858),
@("S08000005",2007,1519,10512,123,286053),
@("S08000006",2007,2929,16873,820,529889),
@("S08000007",2007,5875,41928,495,1191584),
@("S08000008",2007,1631,10212,102,306701),
@("S08000009",2007,3403,20287,149,558139),
@("S08000010",2007,3990,25176,565,801310),
@("S08000011",2007,115,703,3,19770),
@("S08000012",2007,107,740,10,21880),
@("S08000013",2007,1854,14484,279,391639),
@("S08000014",2007,167,846,11,26350)
)

$r = 2
foreach ($row in $rows) {
    $newSheet.Cells.Item($r,1).Value = $row[0]
    $newSheet.Cells.Item($r,2).Value = $row[1]
    $newSheet.Cells.Item($r,3).Value = $row[2]
    $newSheet.Cells.Item($r,4).Value = $row[3]
    $newSheet.Cells.Item($r,5).Value = $row[4]
    $newSheet.Cells.Item($r,6).Value = $row[5]
    $r++
}

# --- 5. Sheet view tweaks ---
$newSheet.Range("D2").Select()

$factSheet.Range("B2").Select()

$prevSheet = $wb.Worksheets.Item("prevelence")
$prevSheet.Range("D26").Select()

$newSheet.Activate()
Write-Host "done"
